# Remove column N from the alcohol measurement data: shift its values
# into column M (the last remaining data column), then delete column N.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 119
for ($r = 1; $r -le $lastRow; $r++) {
    $nVal = $ws.Cells.Item($r, 14).Value()
    $ws.Cells.Item($r, 13).Value = $nVal
}

$ws.Columns.Item(14).Delete()

$ws.Range("M1").Select()
